# SignOutSystem.xlsx - "update function read file import"
#
# Adds the description text in B2 and a small 3-shape flow diagram
# (Input Data -> arrow -> Lay Token va Local Stored) to Sheet1, and
# moves the active selection to W13, matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content -----------------------------------------------------
$ws.Range("B2").Value = "Mô tả chức năng SignOut Sytem"

# --- Shape 1: "Rectangle 1" (Input Data) ------------------------------
$rect1 = $ws.Shapes.AddShape(1, 96, 60, 96.75, 45)
$rect1.Name = "Rectangle 1"
$rect1.Fill.ForeColor.RGB = 16777215
$rect1.TextFrame.Characters().Text = "Input Data"

# --- Shape 2: "Straight Arrow Connector 3" (Input Data -> Lay Token) --
$conn = $ws.Shapes.AddConnector(1, 144, 105, 0.37503937007874016, 831.75)
$conn.Name = "Straight Arrow Connector 3"
$conn.ConnectorFormat.BeginConnect($rect1, 2)

# --- Shape 3: "Rectangle: Rounded Corners 7" (Lay Token va Local Stored)
$rect3 = $ws.Shapes.AddShape(5, 72, 180.00007874015748, 147, 34.5)
$rect3.Name = "Rectangle: Rounded Corners 7"
$rect3.TextFrame.Characters().Text = "Lấy Token và Local  Stored"

# --- Selection ----------------------------------------------------------
[void]$ws.Range("W13").Select()
